$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 6 columns:
#   A: songImage (Image Cover)
#   B: songMusic (Music Title - Filename)
#   C: songName (Music Title)
#   D: artistName (Artist Name)
#   E: albumName (Album Name)
#   F: lyric
#
# The edit simplifies the source sheet down to 3 columns, dropping the
# image/artist/album columns entirely and keeping music/name/lyric
# (old B, C, F) which now become the new A, B, C.

# Overwrite A1:C3 in place with the surviving columns' data (old B/C/F),
# rather than literally deleting columns, so the still-empty placeholder
# cells below the table (A4:A13) aren't disturbed.
$ws.Range("A1").Value2 = "songMusic (Music Title - Filename)"
$ws.Range("B1").Value2 = "songName (Music Title)"
$ws.Range("C1").Value2 = "lyric"

$ws.Range("A2").Value2 = "apapun_yang_terjadi"
$ws.Range("B2").Value2 = "Apapun Yang Terjadi"
$ws.Range("C2").Value2 = "_lyric_apapun_yang_terjadi"

$ws.Range("A3").Value2 = "belum_tidur_feat_sal_priadi"
$ws.Range("B3").Value2 = "Belum Tidur (feat. Sal Priadi)"
$ws.Range("C3").Value2 = "_lyric_belum_tidur_feat_sal_priadi"

# The old A2:A3 (songImage column) carried a distinct cell style; the
# surviving data (old B2:B3) was unstyled, so reset the format here.
$ws.Range("A2:A3").Style = "Normal"

# Drop the now-unused old D:F columns (artistName, albumName, and the
# old location of lyric) along with the old image column's width setup.
$ws.Range("D1:F13").Clear()

# Update the selection to match the new state of the sheet.
$ws.Range("A6").Select()
